# Update column F ("dSF") values on Sheet1 to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -2
    5  = -4
    6  = -9
    7  = -2
    8  = -9
    9  = -4
    10 = 5
    11 = -3
    12 = -5
    13 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
